# Update the "Price" (D) and "Volume(1h)" (E) columns on the crypto list
# sheet with freshly scraped values. Cells with a single decimal point
# (e.g. "1.001") look like numbers to Excel, so a leading apostrophe is
# used to force them to stay plain text, exactly like the source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.945.37'
$ws.Range("E2").Value = '  +0.73%  '
$ws.Range("D3").Value = '1.657.61'
$ws.Range("E3").Value = '  +2.78%  '
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").Value = '''308.97'
$ws.Range("E5").Value = '  +0.80%  '
$ws.Range("D6").Value = '''1.000'
$ws.Range("E6").Value = '  -0.15%  '
$ws.Range("D7").Value = '''0.3882'
$ws.Range("E7").Value = '  -0.25%  '
$ws.Range("D8").Value = '''0.3838'
$ws.Range("E8").Value = '  +1.23%  '
$ws.Range("D9").Value = '''51.15'
$ws.Range("E9").Value = '  +5.30%  '
$ws.Range("D10").Value = '''1.355'
$ws.Range("E10").Value = '  +0.67%  '
$ws.Range("D11").Value = '''1.000'
$ws.Range("E11").Value = '  -0.28%  '
$ws.Range("D12").Value = '''0.08452'
$ws.Range("E12").Value = '  +0.57%  '
$ws.Range("D13").Value = '''23.96'
$ws.Range("E13").Value = '  +1.06%  '
$ws.Range("D14").Value = '''7.153'
$ws.Range("D15").Value = '''7.835'
$ws.Range("E15").Value = '  +5.56%  '
$ws.Range("D16").Value = '''0.00001309'
$ws.Range("E16").Value = '  +3.26%  '
$ws.Range("D17").Value = '1.653.70'
$ws.Range("E17").Value = '  +2.55%  '
$ws.Range("D18").Value = '''94.74'
$ws.Range("E18").Value = '  +1.94%  '
$ws.Range("D19").Value = '''0.06999'
$ws.Range("E19").Value = '  +1.34%  '
$ws.Range("D20").Value = '''19.81'
$ws.Range("E20").Value = '  -0.44%  '
$ws.Range("D21").Value = '''6.903'
$ws.Range("E21").Value = '  +1.94%  '
$ws.Range("D22").Value = '''1.000'
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("D23").Value = '''13.60'
$ws.Range("E23").Value = '  +1.80%  '
$ws.Range("D24").Value = '23.949.33'
$ws.Range("E24").Value = '  +0.66%  '
$ws.Range("D25").Value = '''2.483'
$ws.Range("E25").Value = '  +2.26%  '
$ws.Range("D26").Value = '''3.051'
$ws.Range("E26").Value = '  +9.41%  '
$ws.Range("D27").Value = '''22.09'
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("D28").Value = '''152.69'
$ws.Range("E28").Value = '  -2.68%  '
$ws.Range("D29").Value = '''5.465'
$ws.Range("E29").Value = '  +4.27%  '
$ws.Range("D30").Value = '''139.26'
$ws.Range("E30").Value = '  +0.16%  '
$ws.Range("D31").Value = '''7.794'
$ws.Range("E31").Value = '  +1.06%  '
$ws.Range("D32").Value = '''2.505'
$ws.Range("E32").Value = '  +0.78%  '
$ws.Range("D33").Value = '1.836.15'
$ws.Range("E33").Value = '  +2.59%  '
$ws.Range("D34").Value = '''1.026'
$ws.Range("E34").Value = '  +6.71%  '
$ws.Range("D35").Value = '''0.08028'
$ws.Range("E35").Value = '  -0.71%  '
$ws.Range("D36").Value = '''0.02961'
$ws.Range("E36").Value = '  +3.62%  '
$ws.Range("D37").Value = '''10.99'
$ws.Range("E37").Value = '  +6.07%  '
$ws.Range("D38").Value = '''6.688'
$ws.Range("E38").Value = '  +2.66%  '
$ws.Range("E39").Value = '  +1.93%  '
$ws.Range("D40").Value = '''0.09112'
$ws.Range("E40").Value = '  -0.04%  '
$ws.Range("D41").Value = '''0.7548'
$ws.Range("E41").Value = '  +2.04%  '
$ws.Range("D42").Value = '''13.52'
$ws.Range("E42").Value = '  +1.20%  '
$ws.Range("D43").Value = '''1.418'
$ws.Range("E43").Value = '  -0.26%  '
$ws.Range("D44").Value = '''16.28'
$ws.Range("E44").Value = '  +3.27%  '
$ws.Range("D45").Value = '''0.6940'
$ws.Range("E45").Value = '  +2.22%  '
$ws.Range("D46").Value = '''2.461'
$ws.Range("E46").Value = '  +1.47%  '
$ws.Range("D47").Value = '''4.072'
$ws.Range("E47").Value = '  +0.60%  '
$ws.Range("D48").Value = '''0.9993'
$ws.Range("E48").Value = '  -0.20%  '
$ws.Range("D49").Value = '''0.08289'
$ws.Range("E49").Value = '  +0.96%  '
$ws.Range("D50").Value = '''134.31'
$ws.Range("E50").Value = '  +1.60%  '
$ws.Range("D51").Value = '''1.232'
$ws.Range("E51").Value = '  +3.92%  '
